# CI: data input automation
# Updates the list of selected features in column A (rows 3-46)
# to reflect the latest feature-selection run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @(
    "GarageCars",
    "FullBath",
    "GarageType_Detchd",
    "MSZoning_RM",
    "CentralAir",
    "TotalBsmtSF",
    "Fireplaces",
    "Neighborhood_Crawfor",
    "2ndFlrSF",
    "GarageType_Attchd",
    "BsmtFinSF1",
    "BuildingAge",
    "HeatingQC",
    "BsmtExposure_Gd",
    "MSZoning_RL",
    "SaleType_New",
    "LotArea",
    "Neighborhood_OldTown",
    "GarageFinish_Unf",
    "Exterior1st_BrkFace",
    "BsmtFinType1_GLQ",
    "GarageType_No Garage",
    "MasVnrArea",
    "OverallCond",
    "GarageQual",
    "ScreenPorch",
    "HouseStyle_1Story",
    "SaleCondition_Partial",
    "HouseStyle_2Story",
    "Foundation_PConc",
    "RoofMatl_Tar&Grv",
    "MSZoning_FV",
    "LotFrontage",
    "OpenPorchSF",
    "Neighborhood_NoRidge",
    "LotConfig_CulDSac",
    "RoofStyle_Gambrel",
    "Neighborhood_Timber",
    "Neighborhood_StoneBr",
    "SaleCondition_Normal",
    "BsmtCond",
    "BedroomAbvGr",
    "KitchenAbvGr",
    "SaleType_WD"
)

$startRow = 3
for ($i = 0; $i -lt $features.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $features[$i]
}
